$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update raw metric values ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 260621.35499999998
$metrics.Range("B3").Value = 238185.90999999997
$metrics.Range("B4").Value = 90178.310000000012
$metrics.Range("B5").Value = 10346
$metrics.Range("B6").Value = 841257.125
$metrics.Range("B7").Value = 690980.36
$metrics.Range("B8").Value = 254612.87
$metrics.Range("B9").Value = 33822
$metrics.Range("B10").Value = 34942508.845000006
$metrics.Range("B11").Value = 32736973.149999999
$metrics.Range("B12").Value = 12200426.73
$metrics.Range("B13").Value = 1351729

# Restore the selection on the Metrics sheet to match the saved view state
$metrics.Range("E22").Select() | Out-Null

# --- today sheet: the TODAY()-1 formula in A1 recalculates on its own,
#     and B11:B22/E11:E22/F11:F22 are formulas referencing Metrics!B2:B13,
#     so they update automatically once Metrics changes and the workbook
#     recalculates. We just need to restore the saved selection here too. ---
$today = $wb.Worksheets.Item("today")
$today.Range("E3").Select() | Out-Null

$excel.Calculate()
